$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases & Results")

# ------------------------------------------------------------------
# 1) Prime rows 26 and 27 with the same cell formatting as row 25
#    (the last populated test-case row), then restore column C's
#    original "blank row" style (it is intentionally left untouched
#    by the author in both new rows).
# ------------------------------------------------------------------
$ws.Range("B25:K25").Copy()
$ws.Range("B26:K27").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("C28").Copy()
$ws.Range("C26:C27").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Row heights grow to fit the new wrapped text content
$ws.Rows.Item(26).RowHeight = 43.2
$ws.Rows.Item(27).RowHeight = 57.6

# ------------------------------------------------------------------
# 2) Fix a copy/paste mistake on the existing REQ-22 row: the "Actual
#    Result" referenced test case 12 instead of test case 22.
# ------------------------------------------------------------------
$ws.Range("H25").Value = "Follow the same steps as test case 22, no additional steps required"

# ------------------------------------------------------------------
# 3) New test case row 26: REQ-23
# ------------------------------------------------------------------
$ws.Range("B26").Formula = "=B25+1"
$ws.Range("D26").Value = "REQ-23"
$ws.Range("E26").Value = "Low Impact"
$ws.Range("F26").Value = "Test that LED is turned on for 3 seconds after humidity is high"
$ws.Range("G26").Value = "Humidity must be high at REQ-21"
$ws.Range("H26").Value = "Follow the same steps as test case 22, no additional steps required"
$ws.Range("I26").Value = "LED is turned on for 3 seconds "
$ws.Range("J26").Value = "LED is turned on for 3 seconds "
$ws.Range("K26").Value = "Not Tested"

# ------------------------------------------------------------------
# 4) New test case row 27: REQ-24
# ------------------------------------------------------------------
$ws.Range("B27").Formula = "=B26+1"
$ws.Range("D27").Value = "REQ-24"
$ws.Range("E27").Value = "High Impact"
$ws.Range("F27").Value = "Test that if the return date is late or not by checking from firebase data"
$ws.Range("G27").Value = "The user must scan in REQ-20"
$ws.Range("I27").Value = "The return date is correctly calculated in returnbooks function in main.py"
$ws.Range("J27").Value = "The return date is correctly calculated in returnbooks function in main.py"
$ws.Range("K27").Value = "Not Tested"

# ------------------------------------------------------------------
# 5) Extend the "Not Tested"/"Fail" conditional formatting that was
#    applied to K3:K25 so that it also covers the two new rows.
# ------------------------------------------------------------------
$fc = $ws.Range("K3:K25").FormatConditions.Item(1)
$fc.ModifyAppliesToRange($ws.Range("K3:K27"))

# ------------------------------------------------------------------
# 6) Give the new Test Result cells the same "Pass/Fail/Not Tested"
#    dropdown list validation used throughout column K.
# ------------------------------------------------------------------
$ws.Range("K26:K27").Validation.Add(3, 1, 1, "=Enums!`$B`$2:`$B`$4")

# ------------------------------------------------------------------
# 7) Update the sheet's active selection to reflect where editing
#    finished.
# ------------------------------------------------------------------
$ws.Activate()
$ws.Range("G27").Select()
